$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.257.65"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.108.20"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'574.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'178.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.104.64"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "'0.469"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "3.627.99"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "67.205.95"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "'7.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "3.107.23"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'16.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "'491.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").Value = "'7.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").Value = "'83.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").Value = "'10.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'7.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'2.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'0.111"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "0.0₃0943"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'47.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").Value = "'0.949"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("D39").Value = "'0.313"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").Value = "'49.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "'8.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "'2.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.46%  "
$ws.Range("D45").Value = "2.801.55"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "'371.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "'135.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D50").Value = "'25.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.52%  "
$ws.Range("D51").Value = "'2.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.00%  "
